$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 245; this shifts the existing rows
# 245..284 down to 246..285 and extends the used range to R285.
$ws.Rows.Item(245).Insert()

# Populate the newly inserted row 245 with the new price observation.
$ws.Range("A245").Value = 10
$ws.Range("B245").Value = "Vega Modelo de Temuco"
$ws.Range("C245").Value = "La Araucanía"
$ws.Range("D245").Value = 45034
$ws.Range("E245").Value = 9
$ws.Range("F245").Value = 100112005
$ws.Range("G245").Value = "Puerro"
$ws.Range("H245").Value = "Azul de Maquehue"
$ws.Range("I245").Value = "Primera"
$ws.Range("J245").Value = 65
$ws.Range("K245").Value = 12000
$ws.Range("L245").Value = 12000
$ws.Range("M245").Value = 12000
$ws.Range("N245").Value = "$/docena de paquetes"
$ws.Range("O245").Value = "Provincia de Cautín"
$ws.Range("P245").Value = 1000
$ws.Range("Q245").Value = 12
$ws.Range("R245").Value = "Hortaliza"
